# Append the 2025-03-27 price row to every sheet in the Solar_Prices workbook.
$wb = $excel.ActiveWorkbook

$newDate = "2025-03-27"

$sheetValues = @{
    "N-Dense"                   = "40"
    "N-Type"                    = "43"
    "N-type Wafer"               = "1.21"
    "Cell Topcon 183mm"          = "0.303"
    "Module Topcon 183mm"        = "0.1"
    "Silver Rear_side"           = "5,443"
    "Silver Busbar front-side"   = "8,149"
    "Silver finger front-side"   = "8,199"
    "USD_CNY"                    = "7.2797"
}

foreach ($sheetName in $sheetValues.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Existing rows are plain text (inlineStr) values even though they look
    # like dates/numbers. Use a leading apostrophe to force text entry
    # (otherwise Excel auto-detects a date/number and applies a new
    # NumberFormat style), then clear any formatting flag (quote-prefix)
    # that the apostrophe entry leaves behind so the new cells stay on the
    # default (unstyled) format, matching the rest of the column.
    $lastRow = $ws.Cells.Item(1, 1).End(4).Row
    $newRow = $lastRow + 1

    $dateCell = $ws.Cells.Item($newRow, 1)
    $valueCell = $ws.Cells.Item($newRow, 2)

    $dateCell.Value = "'" + $newDate
    $valueCell.Value = "'" + $sheetValues[$sheetName]

    $ws.Range($dateCell, $valueCell).ClearFormats()
}
